$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in "săpt. 7" (column H) absence values for rows that were missing it
$ws.Range("H7").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H17").Value = 1

# Update existing "săpt. 7" (column H) values
$ws.Range("H9").Value = 2
$ws.Range("H12").Value = 2
$ws.Range("H19").Value = 2
$ws.Range("H21").Value = 2

# Move the active selection on the frozen pane to H19
$ws.Range("H19").Select()
